# TaskList.xlsx - "Updated config file and Added Documentation"
#
# 1) Replace the hard-coded developer / tester e-mail addresses (which lived
#    in the Developer column I and Tester column O) with generic
#    placeholder addresses from the new config.
# 2) Shrink columns I and O now that the new addresses are much shorter.
# 3) Add a new, empty "documentation" row (row 7) with a single cell (O7)
#    formatted using the built-in Hyperlink style (new underlined /
#    hyperlink-colored font + cell style get registered in styles.xml as a
#    side effect, exactly like Excel does the first time a Hyperlink-style
#    cell appears in a workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) swap out the old personal e-mail addresses for the new generic ones
# Developer column (I) first ...
$ws.Range("I2").Value = "dev@gmail.com"
$ws.Range("I3").Value = "dev1@gmail.com"
$ws.Range("I4").Value = "dev@gmail.com"
$ws.Range("I5").Value = "dev1@gmail.com"
$ws.Range("I6").Value = "dev@gmail.com"

# ... then the Tester column (O)
$ws.Range("O2").Value = "tester@gmail.com"
$ws.Range("O3").Value = "tester1@gmail.com"
$ws.Range("O4").Value = "tester@gmail.com"
$ws.Range("O5").Value = "tester1@gmail.com"
$ws.Range("O6").Value = "tester@gmail.com"

# --- 2) the Developer / Tester e-mail columns no longer need to be as wide
$ws.Columns.Item(9).ColumnWidth = 15.666666666666666
$ws.Columns.Item(15).ColumnWidth = 17.666666666666668

# --- 3) new documentation row: an empty, Hyperlink-styled cell at O7.
# Adding (and then removing) a real hyperlink is the natural COM way to get
# Excel to register the built-in "Hyperlink" cell style / font, while
# leaving the target cell itself empty (no address was ever confirmed).
$ws.Hyperlinks.Add($ws.Range("O7"), "", "", "", "") | Out-Null
$ws.Hyperlinks.Delete()
